# Auto-generated edit script
# Applies the cryptocurrency price/volume refresh described by the commit:
#   "Updated cryptos list on Sun Nov 10 10:40:20 UTC 2024 with GitHub Actions"
#
# Strategy per cell:
#   - Many Price (D) values are digit-only strings (e.g. "1.00", "0.119") that
#     Excel would silently re-interpret as numbers (losing trailing zeros) if
#     assigned straight to .Value. For those we force the cell to Text format
#     ("@") before the write so the literal digits are preserved.
#   - After writing, every touched cell is reset to the built-in "Normal" style
#     so no stray formatting delta is introduced versus the original (unstyled)
#     cells - only the cell VALUE changes, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.436.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.192.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "633.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.237"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "3.191.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.593"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +34.95%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.88%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "3.778.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.74%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000229"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +18.47%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "79.271.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "3.188.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.93%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +32.10%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +17.89%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "3.353.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +8.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "77.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.83%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000119"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.00%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.86%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "523.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  +28.32%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.59%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.407"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "Aave"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "192.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "USDe"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.817"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "  +14.64%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.637"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.91%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied 100 cell updates (35 forced to text to preserve literal formatting)"
